$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Beteckning" (A) and "Area (ha)" (G) values between row 2 and row 3,
# and bump the "Förändrad" (C) date forward by one day on both rows.

$ws.Range("A2").Value = "A 36578-2022"
$ws.Range("A3").Value = "A 36523-2022"

$ws.Range("C2").Value = 46066
$ws.Range("C3").Value = 46066

$ws.Range("G2").Value = 0.3
$ws.Range("G3").Value = 0.2
